$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 1126.5
$ws.Range("I31").Value = 1126.5
$ws.Range("K31").Value = 3379.5
$ws.Range("M31").Value = -3149.5
$ws.Range("H40").Value = 4715.36
$ws.Range("I40").Value = 2176
$ws.Range("J40").Value = 4936.174
$ws.Range("K40").Value = 2176
$ws.Range("L40").Value = 4936.174
$ws.Range("M40").Value = -2001
$ws.Range("N40").Value = -5286.174
$ws.Range("H43").Value = 1240.2
$ws.Range("I43").Value = 1228.75
$ws.Range("J43").Value = 1253.2858
$ws.Range("K43").Value = 1228.75
$ws.Range("L43").Value = 1253.2858
$ws.Range("M43").Value = -1159.75
$ws.Range("N43").Value = -1391.2858
$ws.Range("H74").Value = 7613.6924
$ws.Range("I74").Value = 4503
$ws.Range("K74").Value = 4503
$ws.Range("M74").Value = -3567
$ws.Range("H77").Value = 7613.6924
$ws.Range("I77").Value = 4503
$ws.Range("K77").Value = 22515
$ws.Range("M77").Value = -17835
$ws.Range("H116").Value = 4970.231
$ws.Range("I116").Value = 4083.2856
$ws.Range("J116").Value = 6005
$ws.Range("K116").Value = 4083.2856
$ws.Range("L116").Value = 6005
$ws.Range("M116").Value = -641.2856000000002
$ws.Range("N116").Value = -12889
$ws.Range("H137").Value = 108157.3
$ws.Range("I137").Value = 257025.58
$ws.Range("K137").Value = 771076.74
$ws.Range("M137").Value = -768526.74
$ws.Range("H141").Value = 2275.4546
$ws.Range("I141").Value = 1784.2858
$ws.Range("K141").Value = 5352.857400000001
$ws.Range("M141").Value = -172.8574000000008

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 550
$ws.Range("I5").Value = 700
$ws.Range("J5").Value = 100
$ws.Range("K5").Value = 700
$ws.Range("L5").Value = 100
$ws.Range("M5").Value = -588
$ws.Range("N5").Value = -324
$ws.Range("H32").Value = 4414.247
$ws.Range("I32").Value = 3127.8513
$ws.Range("K32").Value = 3127.8513
$ws.Range("M32").Value = -2840.8513
$ws.Range("H102").Value = 4905196.5
$ws.Range("I102").Value = 5558222.5
$ws.Range("K102").Value = 5558222.5
$ws.Range("M102").Value = -5556600.5
$ws.Range("H110").Value = 1394594.1
$ws.Range("I110").Value = 1548660.1
$ws.Range("J110").Value = 7999.5
$ws.Range("K110").Value = 1548660.1
$ws.Range("L110").Value = 7999.5
$ws.Range("M110").Value = -1546615.1
$ws.Range("N110").Value = -12089.5
$ws.Range("H122").Value = 720505.6
$ws.Range("I122").Value = 2119.56
$ws.Range("K122").Value = 6358.68
$ws.Range("M122").Value = -3908.68

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 550
$ws.Range("I4").Value = 700
$ws.Range("J4").Value = 100
$ws.Range("K4").Value = 700
$ws.Range("L4").Value = 100
$ws.Range("M4").Value = -585
$ws.Range("N4").Value = -330
$ws.Range("H76").Value = 100000
$ws.Range("J76").Value = 100000
$ws.Range("L76").Value = 100000
$ws.Range("N76").Value = -100630
$ws.Range("H79").Value = 100000
$ws.Range("J79").Value = 100000
$ws.Range("L79").Value = 100000
$ws.Range("N79").Value = -102184
$ws.Range("H86").Value = 5563861.5
$ws.Range("I86").Value = 5890971
$ws.Range("J86").Value = 3000
$ws.Range("K86").Value = 5890971
$ws.Range("L86").Value = 3000
$ws.Range("M86").Value = -5889848
$ws.Range("N86").Value = -5246
$ws.Range("H89").Value = 5563861.5
$ws.Range("I89").Value = 5890971
$ws.Range("J89").Value = 3000
$ws.Range("K89").Value = 29454855
$ws.Range("L89").Value = 15000
$ws.Range("M89").Value = -29449239
$ws.Range("N89").Value = -26232
$ws.Range("H105").Value = 4167837.2
$ws.Range("I105").Value = 4465454.5
$ws.Range("K105").Value = 4465454.5
$ws.Range("M105").Value = -4463707.5
$ws.Range("H107").Value = 7937833.5
$ws.Range("I107").Value = 8929952
$ws.Range("K107").Value = 8929952
$ws.Range("M107").Value = -8928032
$ws.Range("H134").Value = 4014.1333
$ws.Range("I134").Value = 1510
$ws.Range("J134").Value = 7288.769
$ws.Range("K134").Value = 4530
$ws.Range("L134").Value = 21866.307
$ws.Range("M134").Value = -1995
$ws.Range("N134").Value = -26936.307

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1678.8235
$ws.Range("I16").Value = 1104.091
$ws.Range("J16").Value = 2732.5
$ws.Range("K16").Value = 1104.091
$ws.Range("L16").Value = 2732.5
$ws.Range("M16").Value = -817.0909999999999
$ws.Range("N16").Value = -3306.5
$ws.Range("H31").Value = 27301.416
$ws.Range("J31").Value = 72939.30499999999
$ws.Range("L31").Value = 72939.30499999999
$ws.Range("N31").Value = -73529.30499999999
$ws.Range("H34").Value = 27301.416
$ws.Range("J34").Value = 72939.30499999999
$ws.Range("L34").Value = 72939.30499999999
$ws.Range("N34").Value = -73343.30499999999
$ws.Range("H58").Value = 1926.6072
$ws.Range("I58").Value = 1558.3684
$ws.Range("J58").Value = 2704
$ws.Range("K58").Value = 1558.3684
$ws.Range("L58").Value = 2704
$ws.Range("M58").Value = -1355.3684
$ws.Range("N58").Value = -3110
$ws.Range("H113").Value = 1678.8235
$ws.Range("I113").Value = 1104.091
$ws.Range("J113").Value = 2732.5
$ws.Range("K113").Value = 1104.091
$ws.Range("L113").Value = 2732.5
$ws.Range("M113").Value = 1065.909
$ws.Range("N113").Value = -7072.5
$ws.Range("H136").Value = 1926.6072
$ws.Range("I136").Value = 1558.3684
$ws.Range("J136").Value = 2704
$ws.Range("K136").Value = 4675.1052
$ws.Range("L136").Value = 8112
$ws.Range("M136").Value = -2125.1052
$ws.Range("N136").Value = -13212
$ws.Range("H139").Value = 59114
$ws.Range("J139").Value = 53819
$ws.Range("L139").Value = 53819
$ws.Range("N139").Value = -64099
$ws.Range("H140").Value = 0
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 0
$ws.Range("L140").ClearContents()
$ws.Range("M140").ClearContents()
$ws.Range("N140").Value = 0

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 45
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()
$ws.Range("H46").Value = 278679.9
$ws.Range("I46").Value = 833832.75
$ws.Range("J46").Value = 1103.5
$ws.Range("K46").Value = 2501498.25
$ws.Range("L46").Value = 3310.5
$ws.Range("M46").Value = -2501407.25
$ws.Range("N46").Value = -3492.5
$ws.Range("H75").Value = 170
$ws.Range("I75").Value = 170
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 510
$ws.Range("L75").Value = 0
$ws.Range("M75").ClearContents()
$ws.Range("N75").Value = 488
$ws.Range("H78").Value = 170
$ws.Range("I78").Value = 170
$ws.Range("J78").Value = 0
$ws.Range("K78").Value = 1530
$ws.Range("L78").Value = 0
$ws.Range("M78").ClearContents()
$ws.Range("N78").Value = 3462
$ws.Range("H80").Value = 2780.5
$ws.Range("J80").Value = 2886.2
$ws.Range("L80").Value = 8658.599999999999
$ws.Range("N80").Value = -10530.6
$ws.Range("H83").Value = 2780.5
$ws.Range("J83").Value = 2886.2
$ws.Range("L83").Value = 25975.8
$ws.Range("N83").Value = -35335.8
$ws.Range("H134").Value = 2405.2
$ws.Range("I134").Value = 2339.111
$ws.Range("K134").Value = 7017.333
$ws.Range("M134").Value = -1947.333
$ws.Range("H138").Value = 3570.5715
$ws.Range("I138").Value = 3499
$ws.Range("J138").Value = 4000
$ws.Range("K138").Value = 10497
$ws.Range("L138").Value = 12000
$ws.Range("M138").Value = -5357
$ws.Range("N138").Value = -22280
$ws.Range("H140").Value = 3422.5715
$ws.Range("I140").Value = 2993.1667
$ws.Range("J140").Value = 5999
$ws.Range("K140").Value = 8979.500100000001
$ws.Range("L140").Value = 17997
$ws.Range("M140").Value = -3799.500100000001
$ws.Range("N140").Value = -28357
$ws.Range("H141").Value = 1624.5
$ws.Range("I141").Value = 1624.5
$ws.Range("K141").Value = 4873.5
$ws.Range("M141").Value = 306.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4767
$ws.Range("I7").Value = 2300.1667
$ws.Range("K7").Value = 2300.1667
$ws.Range("M7").Value = -2188.1667
$ws.Range("H22").Value = 127819.43
$ws.Range("J22").Value = 1000
$ws.Range("L22").Value = 1000
$ws.Range("N22").Value = -1590
$ws.Range("H27").Value = 127819.43
$ws.Range("J27").Value = 1000
$ws.Range("L27").Value = 1000
$ws.Range("N27").Value = -1214
$ws.Range("H46").Value = 4083.9473
$ws.Range("I46").Value = 2482.5
$ws.Range("J46").Value = 4823.077
$ws.Range("K46").Value = 2482.5
$ws.Range("L46").Value = 4823.077
$ws.Range("M46").Value = -2294.5
$ws.Range("N46").Value = -5199.077
$ws.Range("H55").Value = 1547.8975
$ws.Range("I55").Value = 1398.0454
$ws.Range("J55").Value = 1741.8235
$ws.Range("K55").Value = 1398.0454
$ws.Range("L55").Value = 1741.8235
$ws.Range("M55").Value = -1225.0454
$ws.Range("N55").Value = -2087.8235
$ws.Range("H100").Value = 3067.1428
$ws.Range("I100").Value = 2828.1304
$ws.Range("K100").Value = 2828.1304
$ws.Range("M100").Value = -2287.1304
$ws.Range("H126").Value = 4767
$ws.Range("I126").Value = 2300.1667
$ws.Range("K126").Value = 6900.500100000001
$ws.Range("M126").Value = -4430.500100000001
$ws.Range("H132").Value = 5373.079
$ws.Range("I132").Value = 4511.4346
$ws.Range("K132").Value = 13534.3038
$ws.Range("M132").Value = -11004.3038
$ws.Range("H136").Value = 71445.5
$ws.Range("I136").Value = 99246.09
$ws.Range("J136").Value = 6577.4443
$ws.Range("K136").Value = 297738.27
$ws.Range("L136").Value = 19732.3329
$ws.Range("M136").Value = -295188.27
$ws.Range("N136").Value = -24832.3329

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2288.1333
$ws.Range("I136").Value = 1756.3077
$ws.Range("K136").Value = 5268.9231
$ws.Range("M136").Value = -2718.9231
